$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 1.1900000000000002
$ws.Range("G2").Value = 1.5090000000000003
$ws.Range("H2").Value = 1.913
$ws.Range("L2").Value = 1.9920000000000002
$ws.Range("M2").Value = 1.8310000000000002
$ws.Range("N2").Value = 1.7429999999999999
$ws.Range("R2").Value = 1.3539999999999999
$ws.Range("W2").Value = 1.8679999999999999
$ws.Range("Z2").Value = 1.675
$ws.Range("AA2").Value = 2.2570000000000001
$ws.Range("AB2").Value = 2.0939999999999999
$ws.Range("AC2").Value = 2.0879999999999996
$ws.Range("AE2").Value = 1.7410000000000001
$ws.Range("AG2").Value = 2.9720000000000004
$ws.Range("AJ2").Value = 1.681
$ws.Range("AN2").Value = 1.0920000000000001
$ws.Range("AQ2").Value = 0.05
$ws.Range("AS2").Value = 1.5839999999999999
$ws.Range("AV2").Value = 0.84899999999999998
$ws.Range("AX2").Value = 1.292
$ws.Range("AY2").Value = 0.63600000000000012
$ws.Range("BB2").Value = 0.85300000000000009
$ws.Range("BD2").Value = 0.82299999999999995
$ws.Range("BI2").Value = 2.0420000000000003
$ws.Range("BJ2").Value = 2.7960000000000003
$ws.Range("BM2").Value = 2.0220000000000002
$ws.Range("BS2").Value = 2.0019999999999998
$ws.Range("BU2").Value = 1.7100000000000002
$ws.Range("BZ2").Value = 2.4560000000000004
$ws.Range("CE2").Value = 1.6929999999999998
$ws.Range("CF2").Value = 2.8089999999999997
$ws.Range("CK2").Value = 2.4329999999999998
$ws.Range("CQ2").Value = 1.3843605864915467
